$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper used to stamp numeric-looking text ("7", "2", "1", "3", ...) into a
# cell as a genuine text value (matching the workbook's existing shared-string
# typed cells) instead of letting Excel auto-convert it to a number.
function Set-TextValue($range, [string]$text) {
    $helper = $ws.Range("Z1")
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $helper.Clear()
}

# Spring AOP now wraps sort(): the concrete algorithm class names are
# replaced everywhere by the CGLIB/JDK dynamic-proxy class name, and the
# measured "Time Taken" column reflects the new proxied run timings.
Set-TextValue $ws.Range("A2") '$Proxy21'
Set-TextValue $ws.Range("C2") "7"

Set-TextValue $ws.Range("A3") '$Proxy21'
Set-TextValue $ws.Range("C3") "2"

Set-TextValue $ws.Range("A4") '$Proxy21'
Set-TextValue $ws.Range("C4") "1"

Set-TextValue $ws.Range("A5") '$Proxy21'
Set-TextValue $ws.Range("C5") "1"

Set-TextValue $ws.Range("A6") '$Proxy21'
Set-TextValue $ws.Range("C6") "3"

Set-TextValue $ws.Range("A7") '$Proxy21'
Set-TextValue $ws.Range("C7") "1"
